# The deck currently has two duplicate "Functions as objects / Callbacks /
# Event-driven" slides (slides 5 and 6) and two duplicate "Lambda" slides
# (slides 12 and 13). Remove the redundant copies, keeping the originals
# that already carry the fully-expanded speaker notes.

$p = $ppt.ActivePresentation

# Delete the duplicate "Functions as objects / Callbacks" slide (index 5).
$p.Slides.Item(5).Delete()

# After the first deletion the duplicate "Lambda" slide has shifted down
# from index 13 to index 12; remove it too.
$p.Slides.Item(12).Delete()
